# Daily attendance processing - 2025-12-21 15:50:56
# Normalises the "Recorded By" (column G) list in each session row:
# the comma-separated recorder names are re-ordered into a stable,
# case-insensitive alphabetical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cellRef = "G" + $row
    $cell = $ws.Range($cellRef)
    $text = $cell.Text

    if ($text -ne $null -and $text -ne "") {
        $parts = $text -split ", "

        if ($parts.Length -gt 1) {
            # Reverse first so that Sort-Object's stable ordering breaks
            # ties (entries equal under case-insensitive comparison,
            # e.g. "System" vs "system") by flipping their relative order.
            $reversed = $parts[($parts.Length - 1)..0]
            $sorted = $reversed | Sort-Object
            $newText = $sorted -join ", "

            # NOTE: -ne/-eq are case-insensitive in this host, which would
            # wrongly treat "System, system" as unchanged when only the
            # case-order of a tie is flipped. Use .Equals() (ordinal,
            # case-sensitive) to detect a genuine change instead.
            if (-not $text.Equals($newText)) {
                $cell.Value = $newText
            }
        }
    }
}
